# Update automàtic: dades i banners [2026-02-28 23:50]
# Applies refreshed scrape timestamps + minor sensor-reading corrections
# to the resum_diari_meteocat sheet (rows 2-46).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('E2').Value = '2026-02-28 23:48:15'
$ws.Range('H2').NumberFormat = '@'
$ws.Range('H2').Value = '79%'
$ws.Range('E3').Value = '2026-02-28 23:48:17'
$ws.Range('E4').Value = '2026-02-28 23:48:20'
$ws.Range('H4').NumberFormat = '@'
$ws.Range('H4').Value = '85%'
$ws.Range('J4').Value = '1025.0 hPa'
$ws.Range('E5').Value = '2026-02-28 23:48:22'
$ws.Range('N5').Value = '-3.6 °C 23:25 TU'
$ws.Range('O5').Value = '-1.6 °C'
$ws.Range('E6').Value = '2026-02-28 23:48:24'
$ws.Range('J6').Value = '1024.9 hPa'
$ws.Range('E7').Value = '2026-02-28 23:48:27'
$ws.Range('E8').Value = '2026-02-28 23:48:29'
$ws.Range('E9').Value = '2026-02-28 23:48:32'
$ws.Range('O9').Value = '11.6 °C'
$ws.Range('E10').Value = '2026-02-28 23:48:34'
$ws.Range('H10').NumberFormat = '@'
$ws.Range('H10').Value = '87%'
$ws.Range('O10').Value = '10.3 °C'
$ws.Range('E11').Value = '2026-02-28 23:48:36'
$ws.Range('H11').NumberFormat = '@'
$ws.Range('H11').Value = '75%'
$ws.Range('E12').Value = '2026-02-28 23:48:39'
$ws.Range('O12').Value = '10.6 °C'
$ws.Range('E13').Value = '2026-02-28 23:48:41'
$ws.Range('E14').Value = '2026-02-28 23:48:44'
$ws.Range('E15').Value = '2026-02-28 23:48:46'
$ws.Range('E16').Value = '2026-02-28 23:48:48'
$ws.Range('H16').NumberFormat = '@'
$ws.Range('H16').Value = '71%'
$ws.Range('N16').Value = '-4.1 °C 23:20 TU'
$ws.Range('E17').Value = '2026-02-28 23:48:51'
$ws.Range('E18').Value = '2026-02-28 23:48:53'
$ws.Range('H18').NumberFormat = '@'
$ws.Range('H18').Value = '85%'
$ws.Range('N18').Value = '6.0 °C 23:29 TU'
$ws.Range('O18').Value = '11.0 °C'
$ws.Range('E19').Value = '2026-02-28 23:48:56'
$ws.Range('I19').Value = '1.7 mm'
$ws.Range('E20').Value = '2026-02-28 23:48:58'
$ws.Range('N20').Value = '-2.4 °C 23:07 TU'
$ws.Range('E21').Value = '2026-02-28 23:49:01'
$ws.Range('J21').Value = '1024.3 hPa'
$ws.Range('E22').Value = '2026-02-28 23:49:03'
$ws.Range('N22').Value = '-4.0 °C 23:29 TU'
$ws.Range('O22').Value = '-1.9 °C'
$ws.Range('E23').Value = '2026-02-28 23:49:05'
$ws.Range('O23').Value = '-0.7 °C'
$ws.Range('E24').Value = '2026-02-28 23:49:08'
$ws.Range('J24').Value = '1025.4 hPa'
$ws.Range('E25').Value = '2026-02-28 23:49:10'
$ws.Range('H25').NumberFormat = '@'
$ws.Range('H25').Value = '67%'
$ws.Range('N25').Value = '-1.7 °C 23:05 TU'
$ws.Range('E26').Value = '2026-02-28 23:49:13'
$ws.Range('H26').NumberFormat = '@'
$ws.Range('H26').Value = '83%'
$ws.Range('I26').Value = '2.4 mm'
$ws.Range('J26').Value = '1024.7 hPa'
$ws.Range('N26').Value = '2.5 °C 23:29 TU'
$ws.Range('O26').Value = '4.6 °C'
$ws.Range('E27').Value = '2026-02-28 23:49:15'
$ws.Range('H27').NumberFormat = '@'
$ws.Range('H27').Value = '61%'
$ws.Range('N27').Value = '-1.1 °C 23:24 TU'
$ws.Range('O27').Value = '1.5 °C'
$ws.Range('E28').Value = '2026-02-28 23:49:17'
$ws.Range('I28').Value = '0.6 mm'
$ws.Range('E29').Value = '2026-02-28 23:49:20'
$ws.Range('E30').Value = '2026-02-28 23:49:22'
$ws.Range('I30').Value = '0.3 mm'
$ws.Range('E31').Value = '2026-02-28 23:49:25'
$ws.Range('H31').NumberFormat = '@'
$ws.Range('H31').Value = '78%'
$ws.Range('E32').Value = '2026-02-28 23:49:27'
$ws.Range('N32').Value = '3.2 °C 23:24 TU'
$ws.Range('E33').Value = '2026-02-28 23:49:29'
$ws.Range('O33').Value = '6.8 °C'
$ws.Range('E34').Value = '2026-02-28 23:49:32'
$ws.Range('I34').Value = '2.9 mm'
$ws.Range('O34').Value = '1.1 °C'
$ws.Range('E35').Value = '2026-02-28 23:49:34'
$ws.Range('H35').NumberFormat = '@'
$ws.Range('H35').Value = '85%'
$ws.Range('E36').Value = '2026-02-28 23:49:36'
$ws.Range('E37').Value = '2026-02-28 23:49:39'
$ws.Range('H37').NumberFormat = '@'
$ws.Range('H37').Value = '82%'
$ws.Range('I37').Value = '2.8 mm'
$ws.Range('E38').Value = '2026-02-28 23:49:41'
$ws.Range('N38').Value = '8.2 °C 23:16 TU'
$ws.Range('E39').Value = '2026-02-28 23:49:43'
$ws.Range('E40').Value = '2026-02-28 23:49:46'
$ws.Range('E41').Value = '2026-02-28 23:49:48'
$ws.Range('E42').Value = '2026-02-28 23:49:51'
$ws.Range('O42').Value = '10.7 °C'
$ws.Range('E43').Value = '2026-02-28 23:49:53'
$ws.Range('O43').Value = '8.0 °C'
$ws.Range('E44').Value = '2026-02-28 23:49:55'
$ws.Range('E45').Value = '2026-02-28 23:49:58'
$ws.Range('N45').Value = '3.8 °C 23:08 TU'
$ws.Range('E46').Value = '2026-02-28 23:50:00'
$ws.Range('O46').Value = '11.2 °C'
